# Update the shared-string text for the Action Category column (rows 2-4)
# and move the active selection to B3, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Land Protection (e.g. conservation easement and/or property acquisition)"
$ws.Range("A4").Value = "Land Management for Protection"
$ws.Range("A3").Value = "Land Management for Protection, Land Protection (e.g. conservation easement and/or property acquisition)" + [char]8221

$ws.Range("B3").Select()
